$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that needs to move from
# 46060 (2026-02-07) to 46061 (2026-02-08) for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
$usedLastRow = $ws.UsedRange.Rows.Count + $ws.UsedRange.Row
if ($usedLastRow -gt $lastRow) { $lastRow = $usedLastRow }
if ($lastRow -lt 328) { $lastRow = 328 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value = 46061
    }
}
